{"js": "// Sequentially find each distinctive cell/title text in the document body\n// and replace it with its updated value (one-to-one text substitution,\n// no structural changes to paragraphs/tables/runs).\nconst replacements = [\n  [\"2025-02-02 Sunday\", \"2025-02-03 Monday\"],\n  [\"400\u00f72=\", \"960\u00f77=\"],\n  [\"947\u00f77=\", \"911\u00f72=\"],\n  [\"859\u00f72=\", \"532\u00f73=\"],\n  [\"682\u00f79=\", \"880\u00f73=\"],\n  [\"305\u00f75=\", \"904\u00f79=\"],\n  [\"804\u00f78=\", \"843\u00f72=\"],\n  [\"418\u00f74=\", \"679\u00f75=\"],\n  [\"364\u00f77=\", \"264\u00f73=\"],\n  [\"285\u00f73=\", \"135\u00f74=\"],\n  [\"899\u00f79=\", \"250\u00f76=\"],\n  [\"224\u00f78=\", \"899\u00f75=\"],\n  [\"545\u00f79=\", \"343\u00f76=\"],\n  [\"732\u00f76=\", \"315\u00f78=\"],\n  [\"707\u00f78=\", \"909\u00f72=\"],\n  [\"681\u00f76=\", \"544\u00f79=\"],\n  [\"503\u00f74=\", \"643\u00f79=\"],\n  [\"420\u00f75=\", \"121\u00f74=\"],\n  [\"937\u00f75=\", \"861\u00f72=\"],\n  [\"549\u00f77=\", \"638\u00f77=\"],\n  [\"716\u00f75=\", \"516\u00f75=\"],\n  [\"161\u00f72=\", \"513\u00f73=\"],\n  [\"195\u00f77=\", \"411\u00f79=\"],\n  [\"992\u00f76=\", \"349\u00f75=\"],\n  [\"979\u00f76=\", \"253\u00f79=\"],\n  [\"281\u00f78=\", \"320\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Sequentially find each distinctive cell/title text in the document\n# and replace it with its updated value (one-to-one text substitution,\n# no structural changes to paragraphs/tables/runs).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-02 Sunday\", \"2025-02-03 Monday\"),\n    @(\"400\u00f72=\", \"960\u00f77=\"),\n    @(\"947\u00f77=\", \"911\u00f72=\"),\n    @(\"859\u00f72=\", \"532\u00f73=\"),\n    @(\"682\u00f79=\", \"880\u00f73=\"),\n    @(\"305\u00f75=\", \"904\u00f79=\"),\n    @(\"804\u00f78=\", \"843\u00f72=\"),\n    @(\"418\u00f74=\", \"679\u00f75=\"),\n    @(\"364\u00f77=\", \"264\u00f73=\"),\n    @(\"285\u00f73=\", \"135\u00f74=\"),\n    @(\"899\u00f79=\", \"250\u00f76=\"),\n    @(\"224\u00f78=\", \"899\u00f75=\"),\n    @(\"545\u00f79=\", \"343\u00f76=\"),\n    @(\"732\u00f76=\", \"315\u00f78=\"),\n    @(\"707\u00f78=\", \"909\u00f72=\"),\n    @(\"681\u00f76=\", \"544\u00f79=\"),\n    @(\"503\u00f74=\", \"643\u00f79=\"),\n    @(\"420\u00f75=\", \"121\u00f74=\"),\n    @(\"937\u00f75=\", \"861\u00f72=\"),\n    @(\"549\u00f77=\", \"638\u00f77=\"),\n    @(\"716\u00f75=\", \"516\u00f75=\"),\n    @(\"161\u00f72=\", \"513\u00f73=\"),\n    @(\"195\u00f77=\", \"411\u00f79=\"),\n    @(\"992\u00f76=\", \"349\u00f75=\"),\n    @(\"979\u00f76=\", \"253\u00f79=\"),\n    @(\"281\u00f78=\", \"320\u00f77=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue = 1, wdReplaceOne = 1, wdReplaceAll = 2\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
